# This script swaps the employee records currently held in rows 16 and 17
# of "Hoja1": the row for "1143384238 DIANA MARCELA DEL VALLE GONZALEZ"
# (161120 / 4028000) and the row for "1143387020 LUCELYS ZUÑIGA VALDEZ"
# (109408 / 2735200) trade places, per the commit "Elimna EC anteriores y
# se agregan nuevos, se modifica base de datos".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Capture current values for the two rows (columns B..G) before overwriting.
$row16 = @{
    B = $ws.Range("B16").Value2
    C = $ws.Range("C16").Value2
    D = $ws.Range("D16").Value2
    E = $ws.Range("E16").Value2
    F = $ws.Range("F16").Value2
    G = $ws.Range("G16").Value2
}

$row17 = @{
    B = $ws.Range("B17").Value2
    C = $ws.Range("C17").Value2
    D = $ws.Range("D17").Value2
    E = $ws.Range("E17").Value2
    F = $ws.Range("F17").Value2
    G = $ws.Range("G17").Value2
}

# Write row17's original data into row16.
$ws.Range("B16").Value2 = $row17.B
$ws.Range("C16").Value2 = $row17.C
$ws.Range("D16").Value2 = $row17.D
$ws.Range("E16").Value2 = $row17.E
$ws.Range("F16").Value2 = $row17.F
$ws.Range("G16").Value2 = $row17.G

# Write row16's original data into row17.
$ws.Range("B17").Value2 = $row16.B
$ws.Range("C17").Value2 = $row16.C
$ws.Range("D17").Value2 = $row16.D
$ws.Range("E17").Value2 = $row16.E
$ws.Range("F17").Value2 = $row16.F
$ws.Range("G17").Value2 = $row16.G
